$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 15:46"

# Updated province/city rows (reordered ranking + refreshed Covid numbers).
# Row => City, Casos totales (B), Casos activos (C), Recuperados (D), Muertes (E)
$updates = @(
    @{ Row = 24; City = "Salamanca";   B = 533; C = 57; D = 430; E = 46 },
    @{ Row = 25; City = "Cantabria";   B = 510; C = 12; D = 484; E = 14 },
    @{ Row = 26; City = "Valladolid";  B = 501; C = 36; D = 444; E = 21 },
    @{ Row = 27; City = "Caceres";     B = 485; C = 3;  D = 447; E = 35 },
    @{ Row = 28; City = "Burgos";      B = 485; C = 55; D = 402; E = 28 },
    @{ Row = 29; City = "Murcia";      B = 477; C = 4;  D = 467; E = 6  },
    @{ Row = 30; City = "Leon";        B = 438; C = 31; D = 376; E = 31 },
    @{ Row = 31; City = "Tenerife";    B = 409; C = 15; D = 400; E = 21 },
    @{ Row = 32; City = "Guadalajara"; B = 404; C = 71; D = 353; E = 49 },
    @{ Row = 35; City = "Segovia";     B = 300; C = 47; D = 219; E = 34 },
    @{ Row = 36; City = "Cordoba";     B = 291; C = 0;  D = 285; E = 6  },
    @{ Row = 37; City = "Cadiz";       B = 278; C = 4;  D = 270; E = 4  },
    @{ Row = 39; City = "Soria";       B = 263; C = 19; D = 230; E = 14 },
    @{ Row = 40; City = "Badajoz";     B = 257; C = 5;  D = 248; E = 4  },
    @{ Row = 41; City = "Avila";       B = 226; C = 36; D = 169; E = 21 },
    @{ Row = 42; City = "Mallorca";    B = 210; C = 18; D = 194; E = 12 },
    @{ Row = 43; City = "Ourense";     B = 189; C = 25; D = 186; E = 3  },
    @{ Row = 48; City = "Zamora";      B = 103; C = 16; D = 79;  E = 8  },
    @{ Row = 49; City = "Palencia";    B = 91;  C = 11; D = 77;  E = 3  },
    @{ Row = 50; City = "Huelva";      B = 77;  C = 2;  D = 74;  E = 1  }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.City
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
}
